$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. SUMMARY paragraph rewrite
# ---------------------------------------------------------------------------
$oldSummary = "Detail-oriented Data Scientist with expertise in Python, machine learning, and data pipeline automation, leveraging strong skills in statistical modeling, data engineering, and cross-functional collaboration."
$newSummary = "Detail-oriented Data Scientist with strong expertise in Python, machine learning (PyTorch, TensorFlow), and data pipeline automation. Proven ability to deploy real-time analytics solutions and optimize models for performance and interpretability. Experienced in statistical modeling, data engineering with Spark, and transforming raw data into actionable business insights."
$d.Content.Find.Execute($oldSummary, $true, $false, $false, $false, $false, $true, 1, $false, $newSummary, 2) | Out-Null

# ---------------------------------------------------------------------------
# 2. "Languages & Frameworks" skill line rewrite
# ---------------------------------------------------------------------------
$oldLangs = "Python, PyTorch, TensorFlow, SQL, Spark"
$newLangs = "Python (NumPy, Pandas, Scikit-learn, Matplotlib), PyTorch, TensorFlow, SQL, Spark"
$d.Content.Find.Execute($oldLangs, $true, $false, $false, $false, $false, $true, 1, $false, $newLangs, 2) | Out-Null

# ---------------------------------------------------------------------------
# 3. Un-bold the existing PROJECTS bullet points and prefix each with "• "
# ---------------------------------------------------------------------------
$bullet = [char]0x2022

$existingBullets = @(
  "Built a big data pipeline to process 94.5M+ rows of NYC DOT Traffic Speeds data using Apache Spark on Google Cloud Dataproc.",
  "Engineered features like hour, weekday, and speed delta to train a PySpark Random Forest classifier for congestion level prediction.",
  "Deployed the model via Flask REST API, enabling real-time inference integrated with Logstash and visualized results using Kibana.",
  "Developed a Retrieval-Augmented Generation (RAG) pipeline to enhance document-based Q&A with LLMs and private knowledge sources.",
  "Used document chunking and vector embeddings (FAISS) for efficient storage and retrieval.",
  "Integrated sub-question decomposition and hierarchical indexing to refine retrieval and improve context relevance."
)

foreach ($b in $existingBullets) {
  $rng = $d.Content
  $found = $rng.Find.Execute($b, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
  if ($found) {
    $rng.Font.Bold = 0
    $rng.InsertBefore("$bullet ")
  }
}

# ---------------------------------------------------------------------------
# 4. Append the new "Leafio" project (a bold title + three un-bold bullets)
#    at the very end of the document, after the RAG project's last bullet.
# ---------------------------------------------------------------------------
for ($i = 0; $i -lt 4; $i++) {
  $endRange = $d.Content
  $endRange.Collapse(0)
  $endRange.InsertParagraphAfter()
}

$count = $d.Paragraphs.Count
$titleIdx = $count - 3
$b1Idx = $count - 2
$b2Idx = $count - 1
$b3Idx = $count

# -- Title paragraph: "Leafio – Real-Time Plant Leaf Disease Detection App" --
$titlePara = $d.Paragraphs($titleIdx)
$titleRange = $titlePara.Range
$titleRange.ListFormat.RemoveNumbers()
$titlePara.Style = "Normal"
$titlePara.Format.LineSpacingRule = 5
$titlePara.Format.LineSpacing = 13.8
$titlePara.Format.SpaceBefore = 12
$titlePara.Format.SpaceAfter = 0
$titleRange.Text = "Leafio " + [char]0x2013 + " Real-Time Plant Leaf Disease Detection App"
$d.Paragraphs($titleIdx).Range.Font.Bold = 1

# -- Three new bullet points --
$newBullets = @(
  "Led the development of an Android app using CNNs for plant disease detection, integrated with a Flask API for real-time predictions, achieving over 85% accuracy.",
  "Designed and trained a custom CNN model using a dataset of diseased and healthy plant leaves, implementing data augmentation techniques to improve model robustness.",
  "Published a research paper on the comparative study of CNN-based disease detection, contributing to AI-driven advancements in precision agriculture."
)
$bulletIdxs = @($b1Idx, $b2Idx, $b3Idx)

for ($i = 0; $i -lt 3; $i++) {
  $idx = $bulletIdxs[$i]
  $p = $d.Paragraphs($idx)
  $p.Range.Text = "$bullet " + $newBullets[$i]
  $d.Paragraphs($idx).Range.Font.Bold = 0
}

Write-Host "Done. Total paragraphs: $($d.Paragraphs.Count)"
